# Update workbook/sheet title and data for 2022-04-20 commit
# (reflecting counts "through 04-12" instead of "through 04-11")

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet tab itself
$ws.Name = "Through 2022-04-12"

# Update the "April (through 04-11)" label to "April (through 04-12)"
$ws.Range("A5").Value = "April (through 04-12)"

# Update April row (row 5) values
$ws.Range("B5").Value = 9
$ws.Range("E5").Value = 20
$ws.Range("G5").Value = 31
$ws.Range("H5").Value = 40
$ws.Range("I5").Value = 48

# Update Total row (row 6) values
$ws.Range("B6").Value = 75
$ws.Range("E6").Value = 217
$ws.Range("G6").Value = 229
$ws.Range("H6").Value = 463
$ws.Range("I6").Value = 482
